$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 of the first effort table ("Matteo") was empty; fill in a new
# effort entry: date 2019-10-27 (serial 43765), topic "Definitions +  Use
# cases" (a new shared string) and effort 3 hrs. Give the topic cell the
# same look as the other topic cells in that table (fill/border already
# present via the existing style) but centered + wrapped, and grow the row
# to fit the two-line label.
$ws.Range("A8").Value = 43765
$ws.Range("B8").Value = "Definitions +  Use cases"
$ws.Range("B8").HorizontalAlignment = -4108
$ws.Range("B8").VerticalAlignment = -4108
$ws.Range("B8").WrapText = $true
$ws.Range("C8").Value = 3
$ws.Rows.Item(8).RowHeight = 30

# Move the active selection like the author's last click before saving.
$ws.Range("H8").Select()
